$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the content of C58:F58 and C59:F59 while keeping their styles
$ws.Range("C58:F58").ClearContents()
$ws.Range("C59:F59").ClearContents()

# Update the selected cell / range on the active sheet view
$ws.Range("F61").Select()
